$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("Z1").Value = 5
$ws.Range("Z1").Font.Name = "Roboto Slab"
$ws.Range("Z1").Font.Size = 16
$ws.Range("Z1").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").Interior.Color = 15773696
